# Applies the "updated 4.0 files and mdl" edit:
#  - About sheet: bump the "last updated" date in C1
#  - MCF sheet: set every non-zero/non-one capacity-factor input to 1

$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = 45392

$mcf = $wb.Worksheets.Item("MCF")
$mcf.Range("B2").Value = 1
$mcf.Range("B3").Value = 1
$mcf.Range("B4").Value = 1
$mcf.Range("B6").Value = 1
$mcf.Range("B10").Value = 1
$mcf.Range("B11").Value = 1
$mcf.Range("B12").Value = 1
$mcf.Range("B13").Value = 1
$mcf.Range("B14").Value = 1
$mcf.Range("B16").Value = 1
$mcf.Range("B17").Value = 1
$mcf.Range("B18").Value = 1

# Selection on MCF sheet moved to B17 per the saved view state.
$mcf.Activate()
$mcf.Range("B17").Select()
